$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.397.73'
$ws.Range('E2').Value = '  +0.80%  '
$ws.Range('D3').Value = '2.485.32'
$ws.Range('E3').Value = '  +0.97%  '
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '522.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.97'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.53%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.558'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.43%  '
$ws.Range('D9').Value = '2.523.49'
$ws.Range('E9').Value = '  +2.30%  '
$ws.Range('E10').Value = '  -0.76%  '
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('E12').Value = '  -1.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.333'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.90%  '
$ws.Range('D14').Value = '2.937.62'
$ws.Range('E14').Value = '  +1.52%  '
$ws.Range('D15').Value = '58.424.74'
$ws.Range('E15').Value = '  +0.94%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.17'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.70%  '
$ws.Range('E17').Value = '  +0.12%  '
$ws.Range('D18').Value = '2.514.50'
$ws.Range('E18').Value = '  +2.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.68'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '322.32'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.16'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.13'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +7.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.996'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.39%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.41'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.404'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.38%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.994'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('E27').Value = '  +0.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.41'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.56%  '
$ws.Range('D29').Value = '0.0₃0755'
$ws.Range('E29').Value = '  +1.67%  '
$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.21'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.89%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.72'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '167.76'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.15%  '
$ws.Range('E33').Value = '  +0.65%  '
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.993'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.44%  '
$ws.Range('E37').Value = '  -3.55%  '
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.47'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.16'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.780'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.64%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '279.39'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.99%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.50'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.12'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.598'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.77%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '123.17'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.92%  '
$ws.Range('E47').Value = '  +0.85%  '
$ws.Range('E48').Value = '  +3.08%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '17.80'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.12%  '
$ws.Range('E50').Value = '  +1.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '16.94'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.62%  '
